# Update Denmark Division 1 odds data (03-04-2024 22:09 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC (2..29) hold the per-match data that moves between rows;
# column A (1) is just the running row index and stays put.
$firstCol = 2
$lastCol = 29

function Get-RowData($row) {
    $data = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $data[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value = $data[$c]
    }
}

# --- Rows 116 and 117 swap their match data ---
$row116 = Get-RowData 116
$row117 = Get-RowData 117
Set-RowData 116 $row117
Set-RowData 117 $row116

# --- Rows 135-139 rotate their match data one position up,
#     with row 135's original data wrapping around to row 139 ---
$row135 = Get-RowData 135
$row136 = Get-RowData 136
$row137 = Get-RowData 137
$row138 = Get-RowData 138
$row139 = Get-RowData 139

Set-RowData 135 $row136
Set-RowData 136 $row137
Set-RowData 137 $row138
Set-RowData 138 $row139
Set-RowData 139 $row135

# --- Row 147: a few odds refreshed ---
$ws.Range("O147").Value = 3.5
$ws.Range("U147").Value = 1.85
$ws.Range("V147").Value = 2

# --- Row 148: a few odds refreshed ---
$ws.Range("R148").Value = 1.875
$ws.Range("S148").Value = 1.975
$ws.Range("U148").Value = 2
$ws.Range("V148").Value = 1.85

# --- Row 149: a few odds refreshed ---
$ws.Range("N149").Value = 1.45
$ws.Range("O149").Value = 4.5
$ws.Range("P149").Value = 6
$ws.Range("Q149").Value = -1.25
$ws.Range("R149").Value = 2.025
$ws.Range("S149").Value = 1.825
$ws.Range("U149").Value = 1.975
$ws.Range("V149").Value = 1.875

# --- Row 150: a single odd refreshed ---
$ws.Range("O150").Value = 3.75
